$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15. This shifts the existing rows 15-21 down to
# 16-22 (preserving their data/styles), matching the diff where each row's
# data moved down by one position and a brand-new row appeared at 15 (with
# the final previously-last row ending up at 22).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly entry.
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 45215
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 100114002
$ws.Range("G15").Value = "Camote"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 17000
$ws.Range("M15").Value = 16500
$ws.Range("N15").Value = "`$/malla 18 kilos"
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 917
$ws.Range("Q15").Value = 18
$ws.Range("R15").Value = "Hortaliza"
